$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 440 (existing rows 440-539 shift down to 442-541)
$ws.Range("A440:A441").EntireRow.Insert()

# New row 440
$ws.Cells.Item(440, 1).Value = 10
$ws.Cells.Item(440, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(440, 3).Value = "La Araucanía"
$ws.Cells.Item(440, 4).Value = 45211
$ws.Cells.Item(440, 5).Value = 9
$ws.Cells.Item(440, 6).Value = 100114013
$ws.Cells.Item(440, 7).Value = "Zanahoria"
$ws.Cells.Item(440, 8).Value = "Sin especificar"
$ws.Cells.Item(440, 9).Value = "Primera"
$ws.Cells.Item(440, 10).Value = 150
$ws.Cells.Item(440, 11).Value = 8000
$ws.Cells.Item(440, 12).Value = 8000
$ws.Cells.Item(440, 13).Value = 8000
$ws.Cells.Item(440, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(440, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(440, 16).Value = 320
$ws.Cells.Item(440, 17).Value = 25
$ws.Cells.Item(440, 18).Value = "Hortaliza"

# New row 441
$ws.Cells.Item(441, 1).Value = 10
$ws.Cells.Item(441, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(441, 3).Value = "La Araucanía"
$ws.Cells.Item(441, 4).Value = 45211
$ws.Cells.Item(441, 5).Value = 9
$ws.Cells.Item(441, 6).Value = 100114013
$ws.Cells.Item(441, 7).Value = "Zanahoria"
$ws.Cells.Item(441, 8).Value = "Sin especificar"
$ws.Cells.Item(441, 9).Value = "Segunda"
$ws.Cells.Item(441, 10).Value = 50
$ws.Cells.Item(441, 11).Value = 6000
$ws.Cells.Item(441, 12).Value = 6000
$ws.Cells.Item(441, 13).Value = 6000
$ws.Cells.Item(441, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(441, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(441, 16).Value = 240
$ws.Cells.Item(441, 17).Value = 25
$ws.Cells.Item(441, 18).Value = "Hortaliza"
